$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) that received new vm_pu values for the 380 kV case
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# One row of values per bus (rows 2-25 => bus index 0-23), in the same column order as $cols
$data = @(
    @(1.02, 1.02737254708353, 1.030327954868506, 1.027428060231597, 1.034689124345579, 1.030976120730964, 1.032530785546686, 1.033139285544224, 1.03024781829231, 1.037487877664863, 1.033997097695105),
    @(1.02, 1.028783061092676, 1.031354870127484, 1.028640958637079, 1.03695936837596, 1.031307795876543, 1.033579089267071, 1.033973945271779, 1.031267344645178, 1.039563476793635, 1.035046890127009),
    @(1.02, 1.029694348594015, 1.032017906832686, 1.029424810083475, 1.03842187591042, 1.031520100847797, 1.034255621116603, 1.034511983723498, 1.031925517281157, 1.040899731321766, 1.035724382729366),
    @(1.02, 1.030077123729111, 1.032296305930334, 1.029754113152411, 1.039035196907896, 1.031608803767672, 1.034539611768453, 1.034737690674292, 1.032201851930859, 1.041459900154175, 1.036008776680516),
    @(1.02, 1.030141374209926, 1.032343030470986, 1.029809391365653, 1.039138088301377, 1.031623665212082, 1.034587270378905, 1.034775559587805, 1.032248228700479, 1.041553862573208, 1.036056502971662),
    @(1.02, 1.029699464541084, 1.032021628150423, 1.029429211135862, 1.038430077047342, 1.031521288257629, 1.03425941746945, 1.034515001527903, 1.031929211090273, 1.040907222536397, 1.03572818447347),
    @(1.02, 1.027849533117296, 1.030675307949818, 1.027838170008886, 1.035457730655642, 1.031088692336402, 1.032885438734724, 1.0334217882787, 1.030592691592625, 1.038190764435993, 1.034352254531343),
    @(1.02, 1.024578595164915, 1.028291632329255, 1.025026845076686, 1.030168673576329, 1.030308559249344, 1.030450354759341, 1.031479562551685, 1.028225639175893, 1.033350391973713, 1.031913712457406),
    @(1.02, 1.022390029082764, 1.026694625936086, 1.023147106582645, 1.026605702730362, 1.029776268973884, 1.028817237606626, 1.030173804330136, 1.026639264455772, 1.030085205260845, 1.030278276091014),
    @(1.02, 1.021440368298322, 1.026001165208688, 1.022331772544813, 1.025053602692837, 1.029542841007108, 1.028107695006803, 1.029605736253679, 1.02595029725261, 1.028661766445323, 1.029567725859313),
    @(1.02, 1.021087312861179, 1.025743284600251, 1.022028704888554, 1.024475638458145, 1.029455689264823, 1.02784377311111, 1.02939432417669, 1.025694068936446, 1.028131553481255, 1.029303429164256),
    @(1.02, 1.021163058586881, 1.025798614467266, 1.022093723826187, 1.024599679899961, 1.029474403867268, 1.027900401934155, 1.029439691263814, 1.025749045145744, 1.028245353839836, 1.029360138406725),
    @(1.02, 1.021411190984946, 1.025979854839739, 1.022306725329315, 1.025005857741694, 1.029535646143166, 1.028085886650192, 1.029588269187343, 1.02592912380695, 1.028617969411827, 1.029545886532335),
    @(1.02, 1.02156403230772, 1.026091483290987, 1.022437933721838, 1.025255924412279, 1.029573320281123, 1.0282001211781, 1.029679758937888, 1.026040034255132, 1.028847352152542, 1.029660283286375),
    @(1.02, 1.022453012089584, 1.026740607122472, 1.023201187656796, 1.026708510400446, 1.029791698471463, 1.028864276635266, 1.030211448476521, 1.026684945079944, 1.030179468439995, 1.030325381920471),
    @(1.02, 1.023010104528222, 1.027147260037298, 1.023679578515051, 1.027617154113399, 1.029927890616399, 1.029280238603009, 1.030544245444357, 1.027088925630004, 1.031012471504329, 1.030741934601919),
    @(1.02, 1.023334855037744, 1.027384266410847, 1.023958481684162, 1.028146255501847, 1.030007045535884, 1.029522631842806, 1.030738103432908, 1.027324362656798, 1.03149742563672, 1.030984672067919),
    @(1.02, 1.023445554261123, 1.027465047854225, 1.024053557909469, 1.0283265149674, 1.030033987327627, 1.029605242808246, 1.030804160602918, 1.027404607203486, 1.031662627299865, 1.031067400350406),
    @(1.02, 1.022950353715708, 1.027103649454173, 1.023628265625254, 1.02751975816188, 1.029913307857311, 1.029235633694633, 1.030508566088202, 1.027045602864852, 1.030923193885314, 1.030697266349456),
    @(1.02, 1.021338130785374, 1.025926492353662, 1.022244007720729, 1.024886288761973, 1.029517624169839, 1.028031276187813, 1.029544527954444, 1.02587610388184, 1.028508284742588, 1.029491198516835),
    @(1.02, 1.020322669092817, 1.025184637380647, 1.021372415043879, 1.023222137076326, 1.029266258470485, 1.027271926779519, 1.028936045288116, 1.025138966528167, 1.026981325095379, 1.028730770745282),
    @(1.02, 1.020861157574898, 1.025578074704005, 1.021834584020175, 1.024105145893802, 1.029399758466228, 1.027674675704054, 1.029258838496286, 1.025529912266876, 1.027791626166164, 1.029134091619487),
    @(1.02, 1.02297735312965, 1.027123355778718, 1.02365145212608, 1.027563769963922, 1.029919898054891, 1.029255789444627, 1.030524688847238, 1.027065179169581, 1.030963537446975, 1.030717450722928),
    @(1.02, 1.025425576431099, 1.02890923788948, 1.025754588054932, 1.031542354839387, 1.030512377769512, 1.031081571508595, 1.031983581523528, 1.028839025519172, 1.034608325085146, 1.032545825606846)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $rowNum = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $cellRef = $cols[$c] + $rowNum.ToString()
        $ws.Range($cellRef).Value = $rowValues[$c]
    }
}
